$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1705.931
$ws.Range("I111").Value = 465.66666
$ws.Range("J111").Value = 1849.0385
$ws.Range("K111").Value = 1396.99998
$ws.Range("L111").Value = 5547.1155
$ws.Range("M111").Value = 1670.00002
$ws.Range("N111").Value = -11681.1155

$ws.Range("H125").Value = 2987.75
$ws.Range("I125").Value = 3012
$ws.Range("J125").Value = 2931.1667
$ws.Range("K125").Value = 27108
$ws.Range("L125").Value = 26380.5003
$ws.Range("M125").Value = -24648
$ws.Range("N125").Value = -31300.5003

$ws.Range("H132").Value = 12931.143
$ws.Range("I132").Value = 12751.117
$ws.Range("J132").Value = 13696.25
$ws.Range("K132").Value = 38253.351
$ws.Range("L132").Value = 41088.75
$ws.Range("M132").Value = -35723.351
$ws.Range("N132").Value = -46148.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3228265.2
$ws.Range("I32").Value = 600.7406999999999
$ws.Range("J32").Value = 25015000
$ws.Range("K32").Value = 600.7406999999999
$ws.Range("L32").Value = 25015000
$ws.Range("M32").Value = -313.7406999999999
$ws.Range("N32").Value = -25015574

$ws.Range("H45").Value = 3108.1667
$ws.Range("J45").Value = 3230.8
$ws.Range("L45").Value = 3230.8
$ws.Range("N45").Value = -3984.8

$ws.Range("H61").Value = 3002.111
$ws.Range("I61").Value = 2853.1667
$ws.Range("K61").Value = 2853.1667
$ws.Range("M61").Value = -2641.1667

$ws.Range("H74").Value = 3286.6086
$ws.Range("I74").Value = 2976.0476
$ws.Range("K74").Value = 2976.0476
$ws.Range("M74").Value = -2102.0476

$ws.Range("H77").Value = 3286.6086
$ws.Range("I77").Value = 2976.0476
$ws.Range("K77").Value = 14880.238
$ws.Range("M77").Value = -10512.238

$ws.Range("H95").Value = 12000
$ws.Range("J95").Value = 12000
$ws.Range("L95").Value = 12000
$ws.Range("N95").Value = -17492

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 2679.7273
$ws.Range("I132").Value = 2547.7
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7643.099999999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -5113.099999999999
$ws.Range("N132").Value = -17060

$ws.Range("H135").Value = 49999.5
$ws.Range("J135").Value = 49999.5
$ws.Range("L135").Value = 49999.5
$ws.Range("N135").Value = -60139.5

$ws.Range("H136").Value = 3002.111
$ws.Range("I136").Value = 2853.1667
$ws.Range("K136").Value = 8559.500100000001
$ws.Range("M136").Value = -6009.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3066.6667
$ws.Range("I20").Value = 3066.6667
$ws.Range("K20").Value = 3066.6667
$ws.Range("M20").Value = -2819.6667

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H86").Value = 1827.6
$ws.Range("I86").Value = 1826.25
$ws.Range("J86").Value = 1833
$ws.Range("K86").Value = 1826.25
$ws.Range("L86").Value = 1833
$ws.Range("M86").Value = -703.25
$ws.Range("N86").Value = -4079

$ws.Range("H89").Value = 1827.6
$ws.Range("I89").Value = 1826.25
$ws.Range("J89").Value = 1833
$ws.Range("K89").Value = 9131.25
$ws.Range("L89").Value = 9165
$ws.Range("M89").Value = -3515.25
$ws.Range("N89").Value = -20397

$ws.Range("H134").Value = 4883.7896
$ws.Range("I134").Value = 1324.5
$ws.Range("K134").Value = 3973.5
$ws.Range("M134").Value = -1438.5

$ws.Range("H140").Value = 130250
$ws.Range("J140").Value = 130250
$ws.Range("L140").Value = 130250
$ws.Range("N140").Value = -140610

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4220.3584
$ws.Range("I31").Value = 2988.8948
$ws.Range("J31").Value = 4908.5293
$ws.Range("K31").Value = 2988.8948
$ws.Range("L31").Value = 4908.5293
$ws.Range("M31").Value = -2693.8948
$ws.Range("N31").Value = -5498.5293

$ws.Range("H34").Value = 4220.3584
$ws.Range("I34").Value = 2988.8948
$ws.Range("J34").Value = 4908.5293
$ws.Range("K34").Value = 2988.8948
$ws.Range("L34").Value = 4908.5293
$ws.Range("M34").Value = -2786.8948
$ws.Range("N34").Value = -5312.5293

$ws.Range("H41").Value = 32974.285
$ws.Range("I41").Value = 22009.166
$ws.Range("K41").Value = 22009.166
$ws.Range("M41").Value = -21581.166

$ws.Range("H51").Value = 44162.25
$ws.Range("J51").Value = 58882.5
$ws.Range("L51").Value = 58882.5
$ws.Range("N51").Value = -60354.5

$ws.Range("H59").Value = 34753
$ws.Range("I59").Value = 13500
$ws.Range("J59").Value = 48921.668
$ws.Range("K59").Value = 13500
$ws.Range("L59").Value = 48921.668
$ws.Range("M59").Value = -12355
$ws.Range("N59").Value = -51211.668

$ws.Range("H60").Value = 14277.667

$ws.Range("H61").Value = 44162.25
$ws.Range("J61").Value = 58882.5
$ws.Range("L61").Value = 58882.5
$ws.Range("N61").Value = -59578.5

$ws.Range("H62").Value = 2333.3333

$ws.Range("H65").Value = 2333.3333

$ws.Range("H68").Value = 35748
$ws.Range("J68").Value = 46246.668
$ws.Range("L68").Value = 46246.668
$ws.Range("N68").Value = -47744.668

$ws.Range("H71").Value = 35748
$ws.Range("J71").Value = 46246.668
$ws.Range("L71").Value = 138740.004
$ws.Range("N71").Value = -146228.004

$ws.Range("H74").Value = 75191.25
$ws.Range("J74").Value = 75191.25
$ws.Range("L74").Value = 75191.25
$ws.Range("N74").Value = -76939.25

$ws.Range("H77").Value = 75191.25
$ws.Range("J77").Value = 75191.25
$ws.Range("L77").Value = 225573.75
$ws.Range("N77").Value = -234309.75

$ws.Range("H107").Value = 888.1
$ws.Range("J107").Value = 1165.8334
$ws.Range("L107").Value = 1165.8334
$ws.Range("N107").Value = -5005.8334

$ws.Range("H134").Value = 2001.2273
$ws.Range("J134").Value = 3708.7144
$ws.Range("L134").Value = 11126.1432
$ws.Range("N134").Value = -16196.1432

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 38.16
$ws.Range("I2").Value = 24.7
$ws.Range("J2").Value = 47.133335
$ws.Range("K2").Value = 148.2
$ws.Range("L2").Value = 282.80001
$ws.Range("M2").Value = -35.19999999999999
$ws.Range("N2").Value = -508.80001

$ws.Range("H8").Value = 406.66666
$ws.Range("I8").Value = 406.66666
$ws.Range("K8").Value = 1219.99998
$ws.Range("M8").Value = -1080.99998

$ws.Range("H92").Value = 6504.2
$ws.Range("J92").Value = 7630.25
$ws.Range("L92").Value = 22890.75
$ws.Range("N92").Value = -25386.75

$ws.Range("H113").Value = 1496.5385
$ws.Range("I113").Value = 682.5
$ws.Range("K113").Value = 2047.5
$ws.Range("M113").Value = 122.5

$ws.Range("H129").Value = 2263.5557
$ws.Range("I129").Value = 723.1667
$ws.Range("J129").Value = 5344.3335
$ws.Range("K129").Value = 2169.5001
$ws.Range("L129").Value = 16033.0005
$ws.Range("M129").Value = 2830.4999
$ws.Range("N129").Value = -26033.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3057
$ws.Range("I102").Value = 3057
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3057
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1435
$ws.Range("N102").ClearContents()

$ws.Range("H126").Value = 3142.5715
$ws.Range("I126").Value = 3142.5715
$ws.Range("K126").Value = 9427.7145
$ws.Range("M126").Value = -6957.7145

$ws.Range("H132").Value = 30351.945
$ws.Range("J132").Value = 3832.3333
$ws.Range("L132").Value = 11496.9999
$ws.Range("N132").Value = -16556.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3366.6667
$ws.Range("I40").Value = 3800
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 3800
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -3664
$ws.Range("N40").Value = -2772

$ws.Range("H46").Value = 2305.5
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 1000
$ws.Range("N46").Value = -1376

$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 4385.909
$ws.Range("I132").Value = 2280
$ws.Range("K132").Value = 6840
$ws.Range("M132").Value = -4310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 25600
$ws.Range("I51").Value = 25600
$ws.Range("K51").Value = 25600
$ws.Range("M51").Value = -25090

$ws.Range("H52").Value = 30021
$ws.Range("I52").Value = 15042
$ws.Range("J52").Value = 45000
$ws.Range("K52").Value = 15042
$ws.Range("L52").Value = 45000
$ws.Range("M52").Value = -14816
$ws.Range("N52").Value = -45452

$ws.Range("H122").Value = 3945.3333
$ws.Range("I122").Value = 3159.8
$ws.Range("K122").Value = 9479.400000000001
$ws.Range("M122").Value = -7029.400000000001

$ws.Range("H132").Value = 2284.7856
$ws.Range("I132").Value = 1809.9
$ws.Range("K132").Value = 5429.700000000001
$ws.Range("M132").Value = -2899.700000000001

$ws.Range("H136").Value = 2341.7896
$ws.Range("I136").Value = 1849.625
$ws.Range("K136").Value = 5548.875
$ws.Range("M136").Value = -2998.875
